# censoring_imp.xlsx — "stabilize cleaning files & internal validity analysis
# for number of pawns balance"
#
# The workbook's table cells are formulas that pull cached values from
# external-workbook links (=[1]decomposition_main_te_0_0!B5, etc). The
# underlying regression-output CSVs were regenerated upstream with slightly
# different numbers; this script pokes the new, refreshed numbers into the
# same cells so the table reflects the rerun.
#
# Numeric-looking results (e.g. "989.9") are written with a leading
# apostrophe so Excel stores them as literal text instead of re-parsing them
# as numbers (which would both lose trailing-zero precision like "0.010" and
# introduce binary floating-point noise like 989.899999...). Cells whose text
# already contains non-numeric markers (stars, parentheses) don't need the
# trick since they can't be mistaken for a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Panel 1 (decomposition_main_te_0_0, rows 4-10) ---
$ws.Range("B5").Value2 = "-236.0***"
$ws.Range("C5").Value2 = "-191.7***"
$ws.Range("D5").Value2 = "'-0.63"
$ws.Range("E5").Value2 = "-75.9**"
$ws.Range("F5").Value2 = "-0.064***"
$ws.Range("D6").Value2 = "(3.01)"
$ws.Range("G9").Value2 = "'0.043"
$ws.Range("B10").Value2 = "'989.9"
$ws.Range("D10").Value2 = "'5.96"
$ws.Range("E10").Value2 = "'396.5"
$ws.Range("F10").Value2 = "'0.44"

# --- Panel 2 (decomposition_main_te_0_1, rows 13-19) ---
$ws.Range("B14").Value2 = "-191.2***"
$ws.Range("D14").Value2 = "'1.17"
$ws.Range("E14").Value2 = "'-15.1"
$ws.Range("F14").Value2 = "'0.0083"
$ws.Range("G14").Value2 = "-0.076***"
$ws.Range("D15").Value2 = "(3.45)"
$ws.Range("D18").Value2 = "'0.004"
$ws.Range("G18").Value2 = "'0.023"
$ws.Range("B19").Value2 = "'989.9"
$ws.Range("D19").Value2 = "'5.96"
$ws.Range("E19").Value2 = "'396.5"
$ws.Range("F19").Value2 = "'0.44"

# --- Panel 5 (decomposition_main_te_imppr, rows 41-48) ---
$ws.Range("B41").Value2 = "-264.9***"
$ws.Range("C41").Value2 = "-169.6***"
$ws.Range("D41").Value2 = "'-1.43"
$ws.Range("E41").Value2 = "-127.4***"
$ws.Range("G41").Value2 = "-0.17***"
$ws.Range("B42").Value2 = "(53.8)"
$ws.Range("C42").Value2 = "(37.2)"
$ws.Range("D42").Value2 = "(3.52)"
$ws.Range("E42").Value2 = "(33.1)"
$ws.Range("F42").Value2 = "(0.025)"
$ws.Range("G42").Value2 = "(0.028)"
$ws.Range("B43").Value2 = "'-42.4"
$ws.Range("C43").Value2 = "'-29.1"
$ws.Range("D43").Value2 = "'-2.66"
$ws.Range("E43").Value2 = "'-14.6"
$ws.Range("F43").Value2 = "'-0.017"
$ws.Range("G43").Value2 = "'0.0026"
$ws.Range("B44").Value2 = "(56.9)"
$ws.Range("C44").Value2 = "(41.8)"
$ws.Range("D44").Value2 = "(3.24)"
$ws.Range("E44").Value2 = "(34.9)"
$ws.Range("G44").Value2 = "(0.029)"
$ws.Range("B47").Value2 = "'0.018"
$ws.Range("D47").Value2 = "'0.002"
$ws.Range("E47").Value2 = "'0.010"
$ws.Range("B48").Value2 = "'1034.5"
$ws.Range("C48").Value2 = "'563.4"
$ws.Range("D48").Value2 = "'7.69"
$ws.Range("E48").Value2 = "'471.2"
$ws.Range("F48").Value2 = "'0.52"
